# Add a new "Greece" worksheet, cloned from "Croatia", with updated data.

$wb = $excel.ActiveWorkbook

$croatia = $wb.Worksheets.Item("Croatia")
$croatia.Activate()
$croatia.Cells.Select()

# Copy Croatia sheet to right after itself
$croatia.Copy($null, $croatia)

# The newly created sheet becomes active and is placed right after Croatia
$greece = $wb.ActiveSheet
$greece.Name = "Greece"

# Update the market name and ticket reference on the new sheet
# (B4 is set first so the new shared-string entries land in the same order
# as the target workbook: NGC ticket before the market name)
$greece.Range("B4").Value = "NGC-4119/T3171"
$greece.Range("B2").Value = "Greece Market"

# Set selection to match the committed state
$greece.Range("I18").Select()

$wb.Activate()
